$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Table 1 ("Type"/"ConnectionString"/...): add two new rows
#    ("Title" and "Module") right before the existing first row.
# ------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$refRow = $t1.Rows.Item(1)

$q1 = [char]8220
$q2 = [char]8221

# Inserting with Add(refRow) always places the new row immediately
# above $refRow, pushing the previously-added row further down, so we
# add the "Module" row first and the "Title" row second to end up
# with Title, Module, Type, ... in that final order.
$moduleRow = $t1.Rows.Add($refRow)
$moduleRow.Cells.Item(1).Range.Text = "Module"
$moduleLine1 = "{0}0{1} View" -f $q1, $q2
$moduleLine2 = "{0}1{1} View con filtro" -f $q1, $q2
$moduleLine3 = "{0}2{1} Grafico" -f $q1, $q2
$moduleRow.Cells.Item(2).Range.Text = $moduleLine1 + "`r" + $moduleLine2 + "`r" + $moduleLine3

$titleRow = $t1.Rows.Add($refRow)
$titleRow.Cells.Item(1).Range.Text = "Title"
$titleRow.Cells.Item(2).Range.Text = "Titolo che verrà visualizzato lato client"

# ------------------------------------------------------------------
# 2) Table 2 ("Modulo client da usare"/"Url"/...): remove the first
#    row, since the module parameter moved server-side.
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$t2.Rows.Item(1).Delete()
